$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("attributes")

# Write new attribute rows for the "solverdportal_experiment_counts" entity.
# Cells are populated in the specific order below so that the workbook's
# shared-strings table grows with new unique strings in this exact sequence:
#   record_url, aggregate_url, hyperlink, comments, text
$ws.Range("B27").Value = "record_url"
$ws.Range("B28").Value = "aggregate_url"
$ws.Range("D27").Value = "hyperlink"
$ws.Range("B26").Value = "comments"
$ws.Range("D26").Value = "text"
$ws.Range("D28").Value = "hyperlink"

$ws.Range("A26").Value = "solverdportal_experiment_counts"
$ws.Range("A27").Value = "solverdportal_experiment_counts"
$ws.Range("A28").Value = "solverdportal_experiment_counts"

# Make "attributes" the active/selected sheet & tab, with the view scrolled
# to show the newly added rows.
$ws.Activate()
$ws.Range("E26").Select()
